# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be written as literal text (matching the
    # original inlineStr cells), preventing Excel from auto-coercing
    # numeric-looking strings (e.g. "1.00", "66.229.69") into numbers.
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "66.229.69"
Set-TextValue $ws.Range("E2") "  -1.92%  "
Set-TextValue $ws.Range("D3") "3.277.55"
Set-TextValue $ws.Range("E3") "  -2.24%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "579.17"
Set-TextValue $ws.Range("E5") "  -0.16%  "
Set-TextValue $ws.Range("E6") "  -3.30%  "
Set-TextValue $ws.Range("D7") "0.627"
Set-TextValue $ws.Range("E7") "  +3.81%  "
Set-TextValue $ws.Range("E8") "  +0.01%  "
Set-TextValue $ws.Range("E9") "  -2.98%  "
Set-TextValue $ws.Range("E10") "  +0.60%  "
Set-TextValue $ws.Range("E11") "  -1.63%  "
Set-TextValue $ws.Range("D12") "3.847.02"
Set-TextValue $ws.Range("E12") "  -2.23%  "
Set-TextValue $ws.Range("E13") "  -3.66%  "
Set-TextValue $ws.Range("D14") "66.245.62"
Set-TextValue $ws.Range("E14") "  -2.17%  "
Set-TextValue $ws.Range("D15") "26.31"
Set-TextValue $ws.Range("E15") "  -4.17%  "
Set-TextValue $ws.Range("B16") "ShibaInu"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D16") "0.0000163"
Set-TextValue $ws.Range("E16") "  -2.73%  "
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.260.31"
Set-TextValue $ws.Range("E17") "  -2.37%  "
Set-TextValue $ws.Range("D18") "434.07"
Set-TextValue $ws.Range("E18") "  -3.07%  "
Set-TextValue $ws.Range("D19") "5.51"
Set-TextValue $ws.Range("E19") "  -3.06%  "
Set-TextValue $ws.Range("D20") "13.16"
Set-TextValue $ws.Range("E20") "  -3.75%  "
Set-TextValue $ws.Range("D21") "7.39"
Set-TextValue $ws.Range("E21") "  -4.61%  "
Set-TextValue $ws.Range("D22") "71.73"
Set-TextValue $ws.Range("E22") "  -3.17%  "
Set-TextValue $ws.Range("E23") "  +0.10%  "
Set-TextValue $ws.Range("D24") "3.421.13"
Set-TextValue $ws.Range("E24") "  -1.88%  "
Set-TextValue $ws.Range("D25") "0.505"
Set-TextValue $ws.Range("E25") "  -1.86%  "
Set-TextValue $ws.Range("E26") "  +1.31%  "
Set-TextValue $ws.Range("E27") "  -6.76%  "
Set-TextValue $ws.Range("E28") "  -2.66%  "
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.05%  "
Set-TextValue $ws.Range("D30") "1.94"
Set-TextValue $ws.Range("E30") "  -1.88%  "
Set-TextValue $ws.Range("D31") "22.26"
Set-TextValue $ws.Range("E31") "  -3.28%  "
Set-TextValue $ws.Range("E32") "  +0.03%  "
Set-TextValue $ws.Range("D33") "5.17"
Set-TextValue $ws.Range("E33") "  -3.73%  "
Set-TextValue $ws.Range("D34") "6.58"
Set-TextValue $ws.Range("E34") "  -3.65%  "
Set-TextValue $ws.Range("E35") "  -4.82%  "
Set-TextValue $ws.Range("D36") "157.04"
Set-TextValue $ws.Range("E36") "  -3.05%  "
Set-TextValue $ws.Range("D37") "1.43"
Set-TextValue $ws.Range("E37") "  -5.44%  "
Set-TextValue $ws.Range("D38") "26.56"
Set-TextValue $ws.Range("E38") "  -3.70%  "
Set-TextValue $ws.Range("D39") "1.78"
Set-TextValue $ws.Range("E39") "  -4.09%  "
Set-TextValue $ws.Range("D40") "2.767.66"
Set-TextValue $ws.Range("E40") "  -2.41%  "
Set-TextValue $ws.Range("D41") "0.772"
Set-TextValue $ws.Range("E41") "  -2.67%  "
Set-TextValue $ws.Range("D42") "4.29"
Set-TextValue $ws.Range("E42") "  -4.58%  "
Set-TextValue $ws.Range("D43") "40.26"
Set-TextValue $ws.Range("E43") "  -0.32%  "
Set-TextValue $ws.Range("D44") "6.04"
Set-TextValue $ws.Range("E44") "  -2.93%  "
Set-TextValue $ws.Range("D45") "0.0658"
Set-TextValue $ws.Range("E45") "  -2.47%  "
Set-TextValue $ws.Range("D46") "319.44"
Set-TextValue $ws.Range("E46") "  -2.13%  "
Set-TextValue $ws.Range("E47") "  -4.40%  "
Set-TextValue $ws.Range("E48") "  -6.61%  "
Set-TextValue $ws.Range("E49") "  -2.97%  "
Set-TextValue $ws.Range("E50") "  +2.22%  "
Set-TextValue $ws.Range("E51") "  +0.02%  "